$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 169166.67
$ws.Cells.Item(106, 9).Value = 251500
$ws.Cells.Item(106, 10).Value = 4500
$ws.Cells.Item(106, 11).Value = 251500
$ws.Cells.Item(106, 12).Value = 4500
$ws.Cells.Item(106, 13).Value = -250869
$ws.Cells.Item(106, 14).Value = -5762

$ws.Cells.Item(113, 8).Value = 86449.234
$ws.Cells.Item(113, 9).Value = 112371.62
$ws.Cells.Item(113, 10).Value = 2201.5
$ws.Cells.Item(113, 11).Value = 112371.62
$ws.Cells.Item(113, 12).Value = 2201.5
$ws.Cells.Item(113, 13).Value = -109117.62
$ws.Cells.Item(113, 14).Value = -8709.5

$ws.Cells.Item(116, 8).Value = 3248.24
$ws.Cells.Item(116, 9).Value = 2405.5293
$ws.Cells.Item(116, 11).Value = 2405.5293
$ws.Cells.Item(116, 13).Value = 1036.4707

$ws.Cells.Item(132, 8).Value = 5323544.5
$ws.Cells.Item(132, 9).Value = 4083.279
$ws.Cells.Item(132, 10).Value = 62507750
$ws.Cells.Item(132, 11).Value = 12249.837
$ws.Cells.Item(132, 12).Value = 187523250
$ws.Cells.Item(132, 13).Value = -9719.837
$ws.Cells.Item(132, 14).Value = -187528310


$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 13159604
$ws.Cells.Item(61, 9).Value = 17858798
$ws.Cells.Item(61, 10).Value = 1865
$ws.Cells.Item(61, 11).Value = 17858798
$ws.Cells.Item(61, 12).Value = 1865
$ws.Cells.Item(61, 13).Value = -17858586
$ws.Cells.Item(61, 14).Value = -2289

$ws.Cells.Item(132, 8).Value = 6759292.5
$ws.Cells.Item(132, 9).Value = 9261377
$ws.Cells.Item(132, 10).Value = 3664
$ws.Cells.Item(132, 11).Value = 27784131
$ws.Cells.Item(132, 12).Value = 10992
$ws.Cells.Item(132, 13).Value = -27781601
$ws.Cells.Item(132, 14).Value = -16052

$ws.Cells.Item(136, 8).Value = 13159604
$ws.Cells.Item(136, 9).Value = 17858798
$ws.Cells.Item(136, 10).Value = 1865
$ws.Cells.Item(136, 11).Value = 53576394
$ws.Cells.Item(136, 12).Value = 5595
$ws.Cells.Item(136, 13).Value = -53573844
$ws.Cells.Item(136, 14).Value = -10695


$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2192.6936
$ws.Cells.Item(105, 9).Value = 1219.6316
$ws.Cells.Item(105, 10).Value = 3733.375
$ws.Cells.Item(105, 11).Value = 1219.6316
$ws.Cells.Item(105, 12).Value = 3733.375
$ws.Cells.Item(105, 13).Value = 527.3684000000001
$ws.Cells.Item(105, 14).Value = -7227.375


$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 8552293
$ws.Cells.Item(31, 9).Value = 5540.6333
$ws.Cells.Item(31, 10).Value = 37041468
$ws.Cells.Item(31, 11).Value = 5540.6333
$ws.Cells.Item(31, 12).Value = 37041468
$ws.Cells.Item(31, 13).Value = -5245.6333
$ws.Cells.Item(31, 14).Value = -37042058

$ws.Cells.Item(34, 8).Value = 8552293
$ws.Cells.Item(34, 9).Value = 5540.6333
$ws.Cells.Item(34, 10).Value = 37041468
$ws.Cells.Item(34, 11).Value = 5540.6333
$ws.Cells.Item(34, 12).Value = 37041468
$ws.Cells.Item(34, 13).Value = -5338.6333
$ws.Cells.Item(34, 14).Value = -37041872

$ws.Cells.Item(58, 8).Value = 1915.826
$ws.Cells.Item(58, 9).Value = 883.0769
$ws.Cells.Item(58, 10).Value = 3258.4
$ws.Cells.Item(58, 11).Value = 883.0769
$ws.Cells.Item(58, 12).Value = 3258.4
$ws.Cells.Item(58, 13).Value = -680.0769
$ws.Cells.Item(58, 14).Value = -3664.4

$ws.Cells.Item(99, 8).Value = 1415.5454
$ws.Cells.Item(99, 9).Value = 1435.7
$ws.Cells.Item(99, 10).Value = 1214
$ws.Cells.Item(99, 11).Value = 1435.7
$ws.Cells.Item(99, 12).Value = 1214
$ws.Cells.Item(99, 13).Value = 62.29999999999995
$ws.Cells.Item(99, 14).Value = -4210

$ws.Cells.Item(126, 8).Value = 1415.5454
$ws.Cells.Item(126, 9).Value = 1435.7
$ws.Cells.Item(126, 10).Value = 1214
$ws.Cells.Item(126, 11).Value = 4307.1
$ws.Cells.Item(126, 12).Value = 3642
$ws.Cells.Item(126, 13).Value = -1837.1
$ws.Cells.Item(126, 14).Value = -8582

$ws.Cells.Item(132, 8).Value = 14288156
$ws.Cells.Item(132, 9).Value = 19232966
$ws.Cells.Item(132, 10).Value = 3147.3333
$ws.Cells.Item(132, 11).Value = 57698898
$ws.Cells.Item(132, 12).Value = 9441.999899999999
$ws.Cells.Item(132, 13).Value = -57696368
$ws.Cells.Item(132, 14).Value = -14501.9999

$ws.Cells.Item(134, 8).Value = 530671.25
$ws.Cells.Item(134, 9).Value = 1554.3055
$ws.Cells.Item(134, 10).Value = 2647139
$ws.Cells.Item(134, 11).Value = 4662.916499999999
$ws.Cells.Item(134, 12).Value = 7941417
$ws.Cells.Item(134, 13).Value = -2127.916499999999
$ws.Cells.Item(134, 14).Value = -7946487

$ws.Cells.Item(136, 8).Value = 1915.826
$ws.Cells.Item(136, 9).Value = 883.0769
$ws.Cells.Item(136, 10).Value = 3258.4
$ws.Cells.Item(136, 11).Value = 2649.2307
$ws.Cells.Item(136, 12).Value = 9775.2
$ws.Cells.Item(136, 13).Value = -99.23070000000007
$ws.Cells.Item(136, 14).Value = -14875.2


$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3204.2156
$ws.Cells.Item(132, 9).Value = 2251.054
$ws.Cells.Item(132, 10).Value = 5723.2856
$ws.Cells.Item(132, 11).Value = 6753.162
$ws.Cells.Item(132, 12).Value = 17169.8568
$ws.Cells.Item(132, 13).Value = -4223.162
$ws.Cells.Item(132, 14).Value = -22229.8568


$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 1942.7858
$ws.Cells.Item(82, 9).Value = 1733.3334
$ws.Cells.Item(82, 10).Value = 2099.875
$ws.Cells.Item(82, 11).Value = 1733.3334
$ws.Cells.Item(82, 12).Value = 2099.875
$ws.Cells.Item(82, 13).Value = -1372.3334
$ws.Cells.Item(82, 14).Value = -2821.875

$ws.Cells.Item(85, 8).Value = 1942.7858
$ws.Cells.Item(85, 9).Value = 1733.3334
$ws.Cells.Item(85, 10).Value = 2099.875
$ws.Cells.Item(85, 11).Value = 1733.3334
$ws.Cells.Item(85, 12).Value = 2099.875
$ws.Cells.Item(85, 13).Value = -485.3334
$ws.Cells.Item(85, 14).Value = -4595.875

$ws.Cells.Item(132, 8).Value = 6585041
$ws.Cells.Item(132, 9).Value = 4006.46
$ws.Cells.Item(132, 11).Value = 12019.38
$ws.Cells.Item(132, 13).Value = -9489.380000000001


$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1561.8182
$ws.Cells.Item(81, 9).Value = 1145
$ws.Cells.Item(81, 10).Value = 1800
$ws.Cells.Item(81, 11).Value = 2290
$ws.Cells.Item(81, 12).Value = 3600
$ws.Cells.Item(81, 13).Value = -1229
$ws.Cells.Item(81, 14).Value = -5722

$ws.Cells.Item(84, 8).Value = 1561.8182
$ws.Cells.Item(84, 9).Value = 1145
$ws.Cells.Item(84, 10).Value = 1800
$ws.Cells.Item(84, 11).Value = 11450
$ws.Cells.Item(84, 12).Value = 18000
$ws.Cells.Item(84, 13).Value = -6146
$ws.Cells.Item(84, 14).Value = -28608

$ws.Cells.Item(122, 8).Value = 1655.6666
$ws.Cells.Item(122, 9).Value = 1815.1482
$ws.Cells.Item(122, 10).Value = 1177.2222
$ws.Cells.Item(122, 11).Value = 5445.444600000001
$ws.Cells.Item(122, 12).Value = 3531.6666
$ws.Cells.Item(122, 13).Value = -2995.444600000001
$ws.Cells.Item(122, 14).Value = -8431.6666

$ws.Cells.Item(132, 8).Value = 1282.338
$ws.Cells.Item(132, 9).Value = 1031.9672
$ws.Cells.Item(132, 10).Value = 2809.6
$ws.Cells.Item(132, 11).Value = 3095.9016
$ws.Cells.Item(132, 12).Value = 8428.8
$ws.Cells.Item(132, 13).Value = -565.9016000000001
$ws.Cells.Item(132, 14).Value = -13488.8

$ws.Cells.Item(136, 8).Value = 901.9423
$ws.Cells.Item(136, 9).Value = 704.5854
$ws.Cells.Item(136, 11).Value = 2113.7562
$ws.Cells.Item(136, 13).Value = 436.2437999999997

